# EIA Table 2.2.C monthly update: October 2016 -> November 2016 edition.
# This mirrors the EIA monthly refresh pattern:
#  - Title / rolling-12-months label text bumped from "October" to "November"
#  - A new monthly data row ("November") is inserted into the "Year 2016"
#    block (row 53, right after the existing Jan-Oct 2016 rows), pushing the
#    "Year to Date" and "Rolling 12 Months Ending in ..." blocks down by one row
#  - The Annual Totals / Year to Date / Rolling 12 months figures are
#    refreshed with the newly released data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Update the headline text that references the reporting month
# ------------------------------------------------------------------
$ws.Range("A2").Value = "by Sector, 2006-November 2016 (Thousand Barrels)"

# ------------------------------------------------------------------
# 2. Insert the new monthly data row for the "Year 2016" block
#    (this shifts rows 53-60 down to 54-61, and Excel automatically
#    re-maps the merged cell / dimension references below it)
# ------------------------------------------------------------------
$ws.Rows("53:53").Insert()

# Copy the number formatting/style from the row above (the last
# existing month row of the block) onto the freshly inserted row so
# that it matches the rest of the table exactly.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row with the November data
$ws.Range("A53").Value = "November"
$ws.Range("B53").Value = 1744
$ws.Range("C53").Value = 1198
$ws.Range("D53").Value = 384
$ws.Range("E53").Value = 19
$ws.Range("F53").Value = 143

# ------------------------------------------------------------------
# 3. Refresh the "Annual Totals" figures (now rows 55-57)
# ------------------------------------------------------------------
$ws.Range("B55").Value = 32709
$ws.Range("C55").Value = 18447
$ws.Range("D55").Value = 11410
$ws.Range("E55").Value = 632
$ws.Range("F55").Value = 2220

$ws.Range("B56").Value = 30277
$ws.Range("C56").Value = 17446
$ws.Range("D56").Value = 10184
$ws.Range("E56").Value = 518
$ws.Range("F56").Value = 2128

$ws.Range("B57").Value = 21622
$ws.Range("C57").Value = 14333
$ws.Range("D57").Value = 5281
$ws.Range("E57").Value = 232
$ws.Range("F57").Value = 1776

# ------------------------------------------------------------------
# 4. Update the "Rolling 12 Months Ending in ..." label (now row 58)
# ------------------------------------------------------------------
$ws.Range("A58").Value = "Rolling 12 Months Ending in November"

# ------------------------------------------------------------------
# 5. Refresh the "Rolling 12 Months" figures (now rows 59-60)
# ------------------------------------------------------------------
$ws.Range("B59").Value = 32197
$ws.Range("C59").Value = 18715
$ws.Range("D59").Value = 10633
$ws.Range("E59").Value = 553
$ws.Range("F59").Value = 2297

$ws.Range("B60").Value = 23413
$ws.Range("C60").Value = 15511
$ws.Range("D60").Value = 5726
$ws.Range("E60").Value = "NM"
$ws.Range("F60").Value = 1931

Write-Host "Update complete"
